$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.128.63"
$ws.Range("E2").Value = "  -3.34%  "

$ws.Range("D3").Value = "1.926.20"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.72"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4728"
$ws.Range("E7").Value = "  -5.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4065"
$ws.Range("E8").Value = "  -3.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.02"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08436"
$ws.Range("E10").Value = "  -9.15%  "

$ws.Range("E11").Value = "  -4.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.26"
$ws.Range("E12").Value = "  -2.75%  "

$ws.Range("D13").Value = "1.920.62"
$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("E14").Value = "  -5.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.099"
$ws.Range("E15").Value = "  -5.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.65"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  -4.05%  "

$ws.Range("E19").Value = "  -2.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("E20").Value = "  -5.78%  "

$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.764"
$ws.Range("E22").Value = "  -3.38%  "

$ws.Range("D23").Value = "28.161.07"
$ws.Range("E23").Value = "  -3.31%  "

$ws.Range("E24").Value = "  -4.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").Value = "2.172.00"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.31"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.12"
$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("E29").Value = "  -4.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.736"
$ws.Range("E30").Value = "  -9.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.75"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9731"
$ws.Range("E32").Value = "  -7.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09604"
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.441"
$ws.Range("E34").Value = "  -5.23%  "

$ws.Range("E35").Value = "  -4.69%  "

$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.035"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02315"
$ws.Range("E38").Value = "  -4.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06173"
$ws.Range("E39").Value = "  -3.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.238"
$ws.Range("E40").Value = "  -6.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6177"
$ws.Range("E41").Value = "  -4.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.06"
$ws.Range("E42").Value = "  -3.86%  "

$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("E44").Value = "  -5.04%  "

$ws.Range("E45").Value = "  -5.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5899"
$ws.Range("E46").Value = "  -5.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.78"
$ws.Range("E47").Value = "  -3.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.038"
$ws.Range("E48").Value = "  -7.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.470"
$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06822"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.03"
$ws.Range("E51").Value = "  -2.98%  "

